{"js": "// Replace each two-digit multiplication expression with its new value.\n// Every occurrence in the document is unique, so a plain text search +\n// whole-match replace is unambiguous for each pair.\nconst replacements = [\n  [\"27\u00d750=\", \"71\u00d767=\"],\n  [\"66\u00d726=\", \"91\u00d785=\"],\n  [\"25\u00d756=\", \"78\u00d747=\"],\n  [\"82\u00d750=\", \"95\u00d768=\"],\n  [\"22\u00d770=\", \"88\u00d754=\"],\n  [\"88\u00d737=\", \"65\u00d746=\"],\n  [\"13\u00d724=\", \"53\u00d769=\"],\n  [\"44\u00d775=\", \"64\u00d745=\"],\n  [\"37\u00d780=\", \"34\u00d746=\"],\n  [\"51\u00d789=\", \"62\u00d786=\"],\n  [\"28\u00d792=\", \"11\u00d788=\"],\n  [\"64\u00d724=\", \"11\u00d747=\"],\n  [\"85\u00d791=\", \"41\u00d732=\"],\n  [\"60\u00d794=\", \"61\u00d735=\"],\n  [\"37\u00d768=\", \"95\u00d732=\"],\n  [\"35\u00d737=\", \"42\u00d799=\"],\n  [\"49\u00d712=\", \"86\u00d780=\"],\n  [\"58\u00d752=\", \"88\u00d746=\"],\n  [\"45\u00d718=\", \"11\u00d753=\"],\n  [\"78\u00d746=\", \"94\u00d743=\"],\n  [\"92\u00d750=\", \"38\u00d715=\"],\n  [\"84\u00d790=\", \"71\u00d743=\"],\n  [\"46\u00d746=\", \"46\u00d779=\"],\n  [\"12\u00d719=\", \"33\u00d783=\"],\n  [\"85\u00d790=\", \"55\u00d769=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: ${oldText}`);\n  }\n\n  for (const result of results.items) {\n    result.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression with its new value.\n# Every expression occurs exactly once in the document, so Find/Replace\n# (wdReplaceAll) on the exact literal text is unambiguous for each pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"27\u00d750=\", \"71\u00d767=\"),\n    @(\"66\u00d726=\", \"91\u00d785=\"),\n    @(\"25\u00d756=\", \"78\u00d747=\"),\n    @(\"82\u00d750=\", \"95\u00d768=\"),\n    @(\"22\u00d770=\", \"88\u00d754=\"),\n    @(\"88\u00d737=\", \"65\u00d746=\"),\n    @(\"13\u00d724=\", \"53\u00d769=\"),\n    @(\"44\u00d775=\", \"64\u00d745=\"),\n    @(\"37\u00d780=\", \"34\u00d746=\"),\n    @(\"51\u00d789=\", \"62\u00d786=\"),\n    @(\"28\u00d792=\", \"11\u00d788=\"),\n    @(\"64\u00d724=\", \"11\u00d747=\"),\n    @(\"85\u00d791=\", \"41\u00d732=\"),\n    @(\"60\u00d794=\", \"61\u00d735=\"),\n    @(\"37\u00d768=\", \"95\u00d732=\"),\n    @(\"35\u00d737=\", \"42\u00d799=\"),\n    @(\"49\u00d712=\", \"86\u00d780=\"),\n    @(\"58\u00d752=\", \"88\u00d746=\"),\n    @(\"45\u00d718=\", \"11\u00d753=\"),\n    @(\"78\u00d746=\", \"94\u00d743=\"),\n    @(\"92\u00d750=\", \"38\u00d715=\"),\n    @(\"84\u00d790=\", \"71\u00d743=\"),\n    @(\"46\u00d746=\", \"46\u00d779=\"),\n    @(\"12\u00d719=\", \"33\u00d783=\"),\n    @(\"85\u00d790=\", \"55\u00d769=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
